# A new weekly cherry-price observation was added to the daily feed for
# "Feria Lagunitas de Puerto Montt". It belongs chronologically right
# after the existing row 38 (2021-12-06) entry, so insert a fresh row at
# position 39 - this pushes the former rows 39..129 down to 40..130 -
# and then populate the new row 39 with the new record's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 39; Excel shifts rows 39:129 -> 40:130 and
# carries the existing row formatting (e.g. the date style on column D)
# down into the freshly inserted row.
$ws.Rows.Item(39).Insert()

# Populate the new row 39 with the new cherry price record.
$ws.Range("A39").Value = 4
$ws.Range("B39").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C39").Value = "Los Lagos"
$ws.Range("D39").Value = 45272
$ws.Range("E39").Value = 10
$ws.Range("F39").Value = "Fruta"
$ws.Range("G39").Value = 100103
$ws.Range("H39").Value = "Frutos de hueso (carozo)"
$ws.Range("I39").Value = 100103001
$ws.Range("J39").Value = "Cereza"
$ws.Range("K39").Value = "Rainier"
$ws.Range("L39").Value = "Primera"
$ws.Range("M39").Value = 450
$ws.Range("N39").Value = 15000
$ws.Range("O39").Value = 15000
$ws.Range("P39").Value = 15000
$ws.Range("Q39").Value = "$/bandeja 10 kilos"
$ws.Range("R39").Value = "Provincia de Curicó"
$ws.Range("S39").Value = 1500
$ws.Range("T39").Value = 10
